$wb = $excel.ActiveWorkbook

# --- "cALL STATISTICS" sheet: fully clear a few stray formatted-but-empty cells ---
$wsCalls = $wb.Worksheets.Item("cALL STATISTICS")
$wsCalls.Range("B2:E2").Clear()
$wsCalls.Range("B9:D9").Clear()
$wsCalls.Range("B11:D11").Clear()

# --- "Sheet5" (Monthly Sales Report) sheet: add a new data row (row 11) and
#     extend the Total / Average formulas down to include it ---
$wsSales = $wb.Worksheets.Item("Sheet5")

$wsSales.Range("D11").Value = 100000000
$wsSales.Range("E11").Formula = "=D11*0.02"

$wsSales.Range("D12").Formula = "=SUM(D5:D11)"
$wsSales.Range("E12").Formula = "=SUM(E5:E11)"
$wsSales.Range("D13").Formula = "=AVERAGE(D5:D11)"
$wsSales.Range("E13").Formula = "=AVERAGE(E5:E11)"

# Move the active selection to D11, matching the saved view state
$wsSales.Range("D11").Select()
